$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K: "Na" header + "NR" for every data row (rows 2-15).
# Copy formatting from the adjacent existing cell first so the new
# cells pick up the same style indices Excel itself would assign,
# then set the value.

# Header K1 -> "Na" (same style as the rest of the header row, from J1)
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Na"

# K2 -> "NR" (style copied from J2)
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("K2").Value = "NR"

# K3 -> "NR" (style copied from J3)
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("K3").Value = "NR"

# K4:K15 -> "NR", filled down from K3's formatting
$ws.Range("K3").Copy()
$ws.Range("K4:K15").PasteSpecial(-4122)
$ws.Range("K4:K15").Value = "NR"

$excel.CutCopyMode = $false

# Match the resulting selection left behind in the saved file
$ws.Range("K16").Select() | Out-Null
